$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

$ws.Range("C2").Value = "z111111e1"
$ws.Range("C3").Value = "Z222222e2"
$ws.Range("C4").Value = "z333333e3"
$ws.Range("C5").Value = "z444444e4"
$ws.Range("C6").Value = "z555555e5"
$ws.Range("C7").Value = "z666666e6"
$ws.Range("C8").Value = "z777777e7"
$ws.Range("C9").Value = "z888888e8"
$ws.Range("C10").Value = "z999999e9"
$ws.Range("C11").Value = "z101010101010e10"
